$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 144 (pushes the existing rows 144-184 down to 146-186)
$ws.Rows("144:145").Insert()

# Row 144: new "Primera" quality entry for the latest week
$ws.Cells.Item(144, 1).Value = 11
$ws.Cells.Item(144, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(144, 3).Value = "Bíobío"
$ws.Cells.Item(144, 4).Value = 44551
$ws.Cells.Item(144, 5).Value = 8
$ws.Cells.Item(144, 6).Value = 100112008
$ws.Cells.Item(144, 7).Value = "Coliflor"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 2000
$ws.Cells.Item(144, 11).Value = 600
$ws.Cells.Item(144, 12).Value = 700
$ws.Cells.Item(144, 13).Value = 650
$ws.Cells.Item(144, 14).Value = "$/unidad"
$ws.Cells.Item(144, 15).Value = "Región Metropolitana"
$ws.Cells.Item(144, 16).Value = 650
$ws.Cells.Item(144, 17).Value = 1
$ws.Cells.Item(144, 18).Value = "Hortaliza"

# Row 145: new "Segunda" quality entry for the same latest week
$ws.Cells.Item(145, 1).Value = 11
$ws.Cells.Item(145, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(145, 3).Value = "Bíobío"
$ws.Cells.Item(145, 4).Value = 44551
$ws.Cells.Item(145, 5).Value = 8
$ws.Cells.Item(145, 6).Value = 100112008
$ws.Cells.Item(145, 7).Value = "Coliflor"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Segunda"
$ws.Cells.Item(145, 10).Value = 1000
$ws.Cells.Item(145, 11).Value = 500
$ws.Cells.Item(145, 12).Value = 500
$ws.Cells.Item(145, 13).Value = 500
$ws.Cells.Item(145, 14).Value = "$/unidad"
$ws.Cells.Item(145, 15).Value = "Región Metropolitana"
$ws.Cells.Item(145, 16).Value = 500
$ws.Cells.Item(145, 17).Value = 1
$ws.Cells.Item(145, 18).Value = "Hortaliza"
